$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": refresh scraped_at timestamps (column K, rows 2-41) ---
$snapshot = $wb.Worksheets.Item("snapshot")
$snapshot.Range("K2").Value = "2025-12-11T07:01:07.581046+00:00"
$snapshot.Range("K3").Value = "2025-12-11T07:01:07.581080+00:00"
$snapshot.Range("K4").Value = "2025-12-11T07:01:07.581104+00:00"
$snapshot.Range("K5").Value = "2025-12-11T07:01:10.523267+00:00"
$snapshot.Range("K6").Value = "2025-12-11T07:01:10.523296+00:00"
$snapshot.Range("K7").Value = "2025-12-11T07:01:10.523316+00:00"
$snapshot.Range("K8").Value = "2025-12-11T07:01:13.370327+00:00"
$snapshot.Range("K9").Value = "2025-12-11T07:01:16.119147+00:00"
$snapshot.Range("K10").Value = "2025-12-11T07:01:20.085392+00:00"
$snapshot.Range("K11").Value = "2025-12-11T07:01:22.961125+00:00"
$snapshot.Range("K12").Value = "2025-12-11T07:01:28.121942+00:00"
$snapshot.Range("K13").Value = "2025-12-11T07:01:28.121985+00:00"
$snapshot.Range("K14").Value = "2025-12-11T07:01:30.380560+00:00"
$snapshot.Range("K15").Value = "2025-12-11T07:01:32.620331+00:00"
$snapshot.Range("K16").Value = "2025-12-11T07:01:35.395503+00:00"
$snapshot.Range("K17").Value = "2025-12-11T07:01:38.301510+00:00"
$snapshot.Range("K18").Value = "2025-12-11T07:01:38.301539+00:00"
$snapshot.Range("K19").Value = "2025-12-11T07:01:40.580721+00:00"
$snapshot.Range("K20").Value = "2025-12-11T07:01:40.580749+00:00"
$snapshot.Range("K21").Value = "2025-12-11T07:01:40.580766+00:00"
$snapshot.Range("K22").Value = "2025-12-11T07:01:43.331174+00:00"
$snapshot.Range("K23").Value = "2025-12-11T07:01:43.331203+00:00"
$snapshot.Range("K24").Value = "2025-12-11T07:01:43.331221+00:00"
$snapshot.Range("K25").Value = "2025-12-11T07:01:43.331238+00:00"
$snapshot.Range("K26").Value = "2025-12-11T07:01:43.331254+00:00"
$snapshot.Range("K27").Value = "2025-12-11T07:01:46.095813+00:00"
$snapshot.Range("K28").Value = "2025-12-11T07:01:46.095845+00:00"
$snapshot.Range("K29").Value = "2025-12-11T07:01:46.095864+00:00"
$snapshot.Range("K30").Value = "2025-12-11T07:01:48.921307+00:00"
$snapshot.Range("K31").Value = "2025-12-11T07:01:48.921339+00:00"
$snapshot.Range("K32").Value = "2025-12-11T07:01:51.614288+00:00"
$snapshot.Range("K33").Value = "2025-12-11T07:01:54.333707+00:00"
$snapshot.Range("K34").Value = "2025-12-11T07:01:54.333743+00:00"
$snapshot.Range("K35").Value = "2025-12-11T07:01:54.333763+00:00"
$snapshot.Range("K36").Value = "2025-12-11T07:01:57.144714+00:00"
$snapshot.Range("K37").Value = "2025-12-11T07:01:57.144743+00:00"
$snapshot.Range("K38").Value = "2025-12-11T07:01:59.646886+00:00"
$snapshot.Range("K39").Value = "2025-12-11T07:01:59.646915+00:00"
$snapshot.Range("K40").Value = "2025-12-11T07:02:02.559066+00:00"
$snapshot.Range("K41").Value = "2025-12-11T07:02:02.559093+00:00"

# --- Sheet "new_injured": clear out the previously-reported new injuries (rows 2-3) ---
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Resize(2).Delete()
